$wb = $excel.ActiveWorkbook

# ===== Sheet LP1912 =====
$ws1 = $wb.Worksheets.Item("LP1912")

# Header updates
$ws1.Range("A2").Value = "Última actualización: 13:41:21"
$ws1.Range("A3").Value = "Total filas: 183"

# Swap C38/C39 (11_ETCHEVERRY <-> 15_ABASTO)
$ws1.Cells.Item(38, 3).Value = "15_ABASTO"
$ws1.Cells.Item(39, 3).Value = "11_ETCHEVERRY"

# Swap rows 62/63 (A, C, D columns)
$ws1.Cells.Item(62, 1).Value = "08:38:24"
$ws1.Cells.Item(62, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(62, 4).Value = 39
$ws1.Cells.Item(63, 1).Value = "07:49:32"
$ws1.Cells.Item(63, 3).Value = "14_ABASTO"
$ws1.Cells.Item(63, 4).Value = 88

# Rows 162-188: updated / re-sorted / newly appended schedule rows
$ws1.Cells.Item(162, 1).Value = "13:41:21"
$ws1.Cells.Item(162, 2).Value = "13:42"
$ws1.Cells.Item(162, 3).Value = "81_EL PELIGRO"
$ws1.Cells.Item(162, 4).Value = 1
$ws1.Cells.Item(162, 5).Value = "LP1912"
$ws1.Cells.Item(163, 1).Value = "13:41:21"
$ws1.Cells.Item(163, 2).Value = "13:44"
$ws1.Cells.Item(163, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(163, 4).Value = 3
$ws1.Cells.Item(163, 5).Value = "LP1912"
$ws1.Cells.Item(164, 1).Value = "11:53:44"
$ws1.Cells.Item(164, 2).Value = "13:47"
$ws1.Cells.Item(164, 3).Value = "225_GOMEZ"
$ws1.Cells.Item(164, 4).Value = 114
$ws1.Cells.Item(164, 5).Value = "LP1912"
$ws1.Cells.Item(165, 1).Value = "12:33:02"
$ws1.Cells.Item(165, 2).Value = "13:54"
$ws1.Cells.Item(165, 3).Value = "15_ABASTO"
$ws1.Cells.Item(165, 4).Value = 81
$ws1.Cells.Item(165, 5).Value = "LP1912"
$ws1.Cells.Item(166, 1).Value = "13:14:31"
$ws1.Cells.Item(166, 2).Value = "14:02"
$ws1.Cells.Item(166, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(166, 4).Value = 48
$ws1.Cells.Item(166, 5).Value = "LP1912"
$ws1.Cells.Item(167, 1).Value = "12:33:02"
$ws1.Cells.Item(167, 2).Value = "14:02"
$ws1.Cells.Item(167, 3).Value = "10_OLMOS"
$ws1.Cells.Item(167, 4).Value = 89
$ws1.Cells.Item(167, 5).Value = "LP1912"
$ws1.Cells.Item(168, 1).Value = "12:46:07"
$ws1.Cells.Item(168, 2).Value = "14:02"
$ws1.Cells.Item(168, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(168, 4).Value = 76
$ws1.Cells.Item(168, 5).Value = "LP1912"
$ws1.Cells.Item(169, 1).Value = "13:14:31"
$ws1.Cells.Item(169, 2).Value = "14:05"
$ws1.Cells.Item(169, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(169, 4).Value = 51
$ws1.Cells.Item(169, 5).Value = "LP1912"
$ws1.Cells.Item(170, 1).Value = "13:41:21"
$ws1.Cells.Item(170, 2).Value = "14:06"
$ws1.Cells.Item(170, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(170, 4).Value = 25
$ws1.Cells.Item(170, 5).Value = "LP1912"
$ws1.Cells.Item(171, 1).Value = "12:46:07"
$ws1.Cells.Item(171, 2).Value = "14:08"
$ws1.Cells.Item(171, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(171, 4).Value = 82
$ws1.Cells.Item(171, 5).Value = "LP1912"
$ws1.Cells.Item(172, 1).Value = "12:53:26"
$ws1.Cells.Item(172, 2).Value = "14:09"
$ws1.Cells.Item(172, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(172, 4).Value = 76
$ws1.Cells.Item(172, 5).Value = "LP1912"
$ws1.Cells.Item(173, 1).Value = "13:41:21"
$ws1.Cells.Item(173, 2).Value = "14:14"
$ws1.Cells.Item(173, 3).Value = "15_ABASTO"
$ws1.Cells.Item(173, 4).Value = 33
$ws1.Cells.Item(173, 5).Value = "LP1912"
$ws1.Cells.Item(174, 1).Value = "12:53:26"
$ws1.Cells.Item(174, 2).Value = "14:16"
$ws1.Cells.Item(174, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(174, 4).Value = 83
$ws1.Cells.Item(174, 5).Value = "LP1912"
$ws1.Cells.Item(175, 1).Value = "12:33:02"
$ws1.Cells.Item(175, 2).Value = "14:17"
$ws1.Cells.Item(175, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(175, 4).Value = 104
$ws1.Cells.Item(175, 5).Value = "LP1912"
$ws1.Cells.Item(176, 1).Value = "12:53:26"
$ws1.Cells.Item(176, 2).Value = "14:17"
$ws1.Cells.Item(176, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(176, 4).Value = 84
$ws1.Cells.Item(176, 5).Value = "LP1912"
$ws1.Cells.Item(177, 1).Value = "12:33:02"
$ws1.Cells.Item(177, 2).Value = "14:18"
$ws1.Cells.Item(177, 3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(177, 4).Value = 105
$ws1.Cells.Item(177, 5).Value = "LP1912"
$ws1.Cells.Item(178, 1).Value = "12:53:26"
$ws1.Cells.Item(178, 2).Value = "14:27"
$ws1.Cells.Item(178, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(178, 4).Value = 94
$ws1.Cells.Item(178, 5).Value = "LP1912"
$ws1.Cells.Item(179, 1).Value = "12:33:02"
$ws1.Cells.Item(179, 2).Value = "14:32"
$ws1.Cells.Item(179, 3).Value = "14X44_ABASTO"
$ws1.Cells.Item(179, 4).Value = 119
$ws1.Cells.Item(179, 5).Value = "LP1912"
$ws1.Cells.Item(180, 1).Value = "12:46:07"
$ws1.Cells.Item(180, 2).Value = "14:34"
$ws1.Cells.Item(180, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(180, 4).Value = 108
$ws1.Cells.Item(180, 5).Value = "LP1912"
$ws1.Cells.Item(181, 1).Value = "12:46:07"
$ws1.Cells.Item(181, 2).Value = "14:39"
$ws1.Cells.Item(181, 3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(181, 4).Value = 113
$ws1.Cells.Item(181, 5).Value = "LP1912"
$ws1.Cells.Item(182, 1).Value = "12:53:26"
$ws1.Cells.Item(182, 2).Value = "14:47"
$ws1.Cells.Item(182, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(182, 4).Value = 114
$ws1.Cells.Item(182, 5).Value = "LP1912"
$ws1.Cells.Item(183, 1).Value = "13:41:21"
$ws1.Cells.Item(183, 2).Value = "14:51"
$ws1.Cells.Item(183, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(183, 4).Value = 70
$ws1.Cells.Item(183, 5).Value = "LP1912"
$ws1.Cells.Item(184, 1).Value = "13:14:31"
$ws1.Cells.Item(184, 2).Value = "14:54"
$ws1.Cells.Item(184, 3).Value = "215_EL PELIGRO"
$ws1.Cells.Item(184, 4).Value = 100
$ws1.Cells.Item(184, 5).Value = "LP1912"
$ws1.Cells.Item(185, 1).Value = "13:14:31"
$ws1.Cells.Item(185, 2).Value = "15:02"
$ws1.Cells.Item(185, 3).Value = "10_OLMOS"
$ws1.Cells.Item(185, 4).Value = 108
$ws1.Cells.Item(185, 5).Value = "LP1912"
$ws1.Cells.Item(186, 1).Value = "13:14:31"
$ws1.Cells.Item(186, 2).Value = "15:13"
$ws1.Cells.Item(186, 3).Value = "17X38_ROMERO"
$ws1.Cells.Item(186, 4).Value = 119
$ws1.Cells.Item(186, 5).Value = "LP1912"
$ws1.Cells.Item(187, 1).Value = "13:41:21"
$ws1.Cells.Item(187, 2).Value = "15:18"
$ws1.Cells.Item(187, 3).Value = "14_ABASTO"
$ws1.Cells.Item(187, 4).Value = 97
$ws1.Cells.Item(187, 5).Value = "LP1912"
$ws1.Cells.Item(188, 1).Value = "13:41:21"
$ws1.Cells.Item(188, 2).Value = "15:34"
$ws1.Cells.Item(188, 3).Value = "215C_EL PATO"
$ws1.Cells.Item(188, 4).Value = 113
$ws1.Cells.Item(188, 5).Value = "LP1912"

# ===== Sheet LP1912-215 =====
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 13:41:21"
$ws2.Range("A3").Value = "Total filas: 32"
$ws2.Cells.Item(37, 1).Value = "13:41:21"
$ws2.Cells.Item(37, 2).Value = "15:34"
$ws2.Cells.Item(37, 3).Value = "215C_EL PATO"
$ws2.Cells.Item(37, 4).Value = 113
$ws2.Cells.Item(37, 5).Value = "LP1912"

# ===== Sheet 6203-6173 =====
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 13:41:21"
$ws3.Range("A3").Value = "Total filas: 28"
$ws3.Cells.Item(33, 1).Value = "13:41:21"
$ws3.Cells.Item(33, 2).Value = "15:22"
$ws3.Cells.Item(33, 3).Value = "215A_LA PLATA"
$ws3.Cells.Item(33, 4).Value = 101
$ws3.Cells.Item(33, 5).Value = "L6173"
